$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("unet")
$ws2 = $wb.Worksheets.Item("simple cnn")

# --- numeric-only fills first (no new shared strings) ---

# Sheet1 row 38: fill in remaining numeric data for existing "Notes" label in A38
$ws1.Range("B38").Value = 0.1368
$ws1.Range("C38").Value = 0.5165

# Sheet1 row 39: new row - numeric cells
$ws1.Range("B39").Value = 0.0891
$ws1.Range("C39").Value = 0.6173

# Sheet2 row 19: fill in remaining numeric data for existing label in A19
$ws2.Range("B19").Value = 0.234
$ws2.Range("C19").Value = -0.6168
$ws2.Range("D19").Value = 0.9041
$ws2.Range("E19").Value = 0.0945

# Sheet2 row 20: new row - numeric cells
$ws2.Range("B20").Value = 0.186
$ws2.Range("C20").Value = -0.6996
$ws2.Range("D20").Value = 0.9271
$ws2.Range("E20").Value = 0.0722

# Sheet2 row 21: new row

# --- new shared-string text, in the order originally typed ---

# 1) sheet2 A20 augmentor description
$ws2.Range("A20").Value = "iaa.SomeOf((0, 2), [
                iaa.Flipud(1.0),
                iaa.Multiply((0.5, 1.5), per_channel=0.5),
                iaa.Affine(shear=(-20, 20)),
                iaa.Affine(translate_px={""x"": (-15, 15), ""y"": (-15, 15)})
                ])"
$ws2.Range("A20").WrapText = $true
$ws2.Rows.Item(20).RowHeight = 114

# 2) sheet1 D38 note
$ws1.Range("D38").Value = "early stop at 30"

# 3) sheet1 D37 note (added after D38)
$ws1.Range("D37").Value = "?early stop at 35ish"

# 4) sheet1 A39 augmentor description
$ws1.Range("A39").Value = "iaa.SomeOf((0, 2), [
        iaa.Fliplr(1.0),  # horizontal flip
        iaa.Flipud(1.0),  # vertical flip
        iaa.Affine(translate_px={""x"": (-15, 15), ""y"": (-15, 15)}),  # lighten or darken
    ])"
$ws1.Range("A39").WrapText = $true
$ws1.Rows.Item(39).RowHeight = 85.5

# 5) sheet1 D39 note
$ws1.Range("D39").Value = "early stop at 39"

# 6) sheet2 A21 augmentor description
$ws2.Range("A21").Value = "iaa.SomeOf((0, 2), [
                iaa.Flipud(1.0),
                iaa.Fliplr(1.0),
                iaa.Multiply((0.5, 1.5), per_channel=0.5),
                iaa.Affine(shear=(-20, 20)),
                iaa.Affine(rotate=(-15, 15)),
                ])"
$ws2.Range("A21").WrapText = $true
$ws2.Rows.Item(21).RowHeight = 114

# --- final view/selection state ---
$ws1.Range("B44").Select()
$ws2.Range("A29").Select()
